$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B..AB) must be swapped with each other.
# Column A (the running id/index) is left untouched.
$pairs = @(
    @(140, 141),
    @(260, 261),
    @(293, 294)
)

$firstCol = 2   # column B
$lastCol  = 28  # column AB

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
